$wb = $excel.ActiveWorkbook

# The file "b2545090-b8fe-4df0-8160-19d983cec574.md" (row 3 on every sheet)
# is now ready for handoff again, with a fresh handoff timestamp recorded
# per-locale on the zh-cn and de-de detail sheets.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-18 02:16:43"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-18 02:16:55"
